# Updates the EPEX Spot prices workbook with one additional day of data:
#   - "Prix Spot" sheet: new column BX ("28-aug" header + 24 hourly prices)
#   - "Gaz" sheet: new row 73 (2025-08-26 + price)
#   - "CO2" sheet: new row 73 (2025-08-26 + price)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Prix Spot" -> add column BX
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# BX1 header cell: copy format (bold, centered, bordered) from BW1, then set text
$ws1.Range("BW1").Copy($ws1.Range("BX1"))
$ws1.Range("BX1").Value = "28-aug"

# BX2:BX25 hourly values (plain numbers, default formatting like the rest of the data)
$ws1.Range("BX2").Value = 97.15000000000001
$ws1.Range("BX3").Value = 90.14
$ws1.Range("BX4").Value = 86.06999999999999
$ws1.Range("BX5").Value = 73.25
$ws1.Range("BX6").Value = 69.40000000000001
$ws1.Range("BX7").Value = 76.42
$ws1.Range("BX8").Value = 83.38
$ws1.Range("BX9").Value = 90.38
$ws1.Range("BX10").Value = 90.38
$ws1.Range("BX11").Value = 69.40000000000001
$ws1.Range("BX12").Value = 42.57
$ws1.Range("BX13").Value = 20.46
$ws1.Range("BX14").Value = 17.21
$ws1.Range("BX15").Value = 2.09
$ws1.Range("BX16").Value = 0.65
$ws1.Range("BX17").Value = 0.01
$ws1.Range("BX18").Value = 0.65
$ws1.Range("BX19").Value = 21.04
$ws1.Range("BX20").Value = 42.46
$ws1.Range("BX21").Value = 71.12
$ws1.Range("BX22").Value = 90
$ws1.Range("BX23").Value = 97.08
$ws1.Range("BX24").Value = 93.77
$ws1.Range("BX25").Value = 88.09

# ---------------------------------------------------------------------------
# Sheet 2: "Gaz" -> add row 73
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Column A holds a date-formatted text string (e.g. "2025-08-25") stored as plain
# text, not a real date. Force text formatting, assign, then restore the default
# "Normal" style so the cell ends up unstyled just like the existing rows.
$a73 = $ws2.Range("A73")
$a73.NumberFormat = "@"
$a73.Value = "2025-08-26"
$a73.Style = "Normal"

$ws2.Range("B73").Value = 32.175

# ---------------------------------------------------------------------------
# Sheet 3: "CO2" -> add row 73
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$a73b = $ws3.Range("A73")
$a73b.NumberFormat = "@"
$a73b.Value = "2025-08-26"
$a73b.Style = "Normal"

$ws3.Range("B73").Value = 72.28
